# Add a "last updated" style date stamp to the About sheet (cell C1),
# formatted as a short date (built-in numFmtId 14).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

$ws.Range("C1").Value = 44307
$ws.Range("C1").NumberFormat = "mm-dd-yy"
